# Update the "Förändrad" date column (C) for all data rows (2-216)
# from serial date 45179 (2023-09-10) to 45180 (2023-09-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 216
$col = 3  # Column C

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    if ($cell.Value -eq 45179) {
        $cell.Value = 45180
    }
}
